# Auto-generated by build script: updates Leveling profit-tracking values
# across all 8 worksheets to match the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 6148.1665
$ws.Range("J32").Value = 4347.75
$ws.Range("L32").Value = 4347.75
$ws.Range("N32").Value = -4999.75
# Row 80
$ws.Range("H80").Value = 423.21054
$ws.Range("I80").Value = 174.875
$ws.Range("J80").Value = 603.8182
$ws.Range("K80").Value = 524.625
$ws.Range("L80").Value = 1811.4546
$ws.Range("M80").Value = 473.375
$ws.Range("N80").Value = -3807.4546
# Row 83
$ws.Range("H83").Value = 423.21054
$ws.Range("I83").Value = 174.875
$ws.Range("J83").Value = 603.8182
$ws.Range("K83").Value = 1573.875
$ws.Range("L83").Value = 5434.3638
$ws.Range("M83").Value = 3418.125
$ws.Range("N83").Value = -15418.3638
# Row 132
$ws.Range("H132").Value = 1071.1803
$ws.Range("I132").Value = 707.8596
$ws.Range("J132").Value = 6248.5
$ws.Range("K132").Value = 2123.5788
$ws.Range("L132").Value = 18745.5
$ws.Range("M132").Value = 406.4211999999998
$ws.Range("N132").Value = -23805.5

$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 286571.28
$ws.Range("I4").Value = 286571.28
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 286571.28
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -286455.28
$ws.Range("N4").ClearContents()
# Row 5
$ws.Range("H5").Value = 205.21053
$ws.Range("I5").Value = 171.36363
$ws.Range("J5").Value = 251.75
$ws.Range("K5").Value = 171.36363
$ws.Range("L5").Value = 251.75
$ws.Range("M5").Value = -59.36363
$ws.Range("N5").Value = -475.75
# Row 28
$ws.Range("H28").Value = 15360.429
$ws.Range("I28").Value = 6621.25
$ws.Range("K28").Value = 6621.25
$ws.Range("M28").Value = -6429.25
# Row 45
$ws.Range("H45").Value = 391489
$ws.Range("I45").Value = 596818.25
$ws.Range("J45").Value = 3644.889
$ws.Range("K45").Value = 596818.25
$ws.Range("L45").Value = 3644.889
$ws.Range("M45").Value = -596441.25
$ws.Range("N45").Value = -4398.889
# Row 99
$ws.Range("H99").Value = 15360.429
$ws.Range("I99").Value = 6621.25
$ws.Range("K99").Value = 6621.25
$ws.Range("M99").Value = -3626.25
# Row 110
$ws.Range("H110").Value = 2931.6191
$ws.Range("I110").Value = 2978.2
$ws.Range("K110").Value = 2978.2
$ws.Range("M110").Value = -933.1999999999998
# Row 132
$ws.Range("H132").Value = 2037
$ws.Range("I132").Value = 1199.75
$ws.Range("J132").Value = 7395.4
$ws.Range("K132").Value = 3599.25
$ws.Range("L132").Value = 22186.2
$ws.Range("M132").Value = -1069.25
$ws.Range("N132").Value = -27246.2

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 205.21053
$ws.Range("I4").Value = 171.36363
$ws.Range("J4").Value = 251.75
$ws.Range("K4").Value = 171.36363
$ws.Range("L4").Value = 251.75
$ws.Range("M4").Value = -56.36363
$ws.Range("N4").Value = -481.75
# Row 22
$ws.Range("H22").Value = 566.6667
$ws.Range("I22").Value = 566.6667
$ws.Range("K22").Value = 566.6667
$ws.Range("M22").Value = -393.6667
# Row 94
$ws.Range("H94").Value = 552.8929000000001
$ws.Range("I94").Value = 552.8929000000001
$ws.Range("K94").Value = 552.8929000000001
$ws.Range("M94").Value = -101.8929000000001
# Row 105
$ws.Range("H105").Value = 5501.5884
$ws.Range("I105").Value = 5934.4
$ws.Range("K105").Value = 5934.4
$ws.Range("M105").Value = -4187.4
# Row 134
$ws.Range("H134").Value = 1687.1428
$ws.Range("I134").Value = 1687.1428
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5061.428400000001
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -2526.428400000001
$ws.Range("N134").ClearContents()
# Row 141
$ws.Range("H141").Value = 92082.5
$ws.Range("J141").Value = 123456
$ws.Range("L141").Value = 123456
$ws.Range("N141").Value = -133816

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 614
$ws.Range("I22").Value = 614
$ws.Range("K22").Value = 614
$ws.Range("M22").Value = -264
# Row 31
$ws.Range("H31").Value = 3263.32
$ws.Range("I31").Value = 1764.3158
$ws.Range("K31").Value = 1764.3158
$ws.Range("M31").Value = -1469.3158
# Row 34
$ws.Range("H34").Value = 3263.32
$ws.Range("I34").Value = 1764.3158
$ws.Range("K34").Value = 1764.3158
$ws.Range("M34").Value = -1562.3158
# Row 132
$ws.Range("H132").Value = 1595
$ws.Range("I132").Value = 1556.7858
$ws.Range("J132").Value = 2130
$ws.Range("K132").Value = 4670.357400000001
$ws.Range("L132").Value = 6390
$ws.Range("M132").Value = -2140.357400000001
$ws.Range("N132").Value = -11450

$ws = $wb.Worksheets.Item("CUL")
# Row 55
$ws.Range("H55").Value = 9618349
$ws.Range("I55").Value = 615.5
$ws.Range("J55").Value = 13892898
$ws.Range("K55").Value = 1846.5
$ws.Range("L55").Value = 41678694
$ws.Range("M55").Value = -1669.5
$ws.Range("N55").Value = -41679048
# Row 92
$ws.Range("H92").Value = 593.5
$ws.Range("J92").Value = 887.5
$ws.Range("L92").Value = 2662.5
$ws.Range("N92").Value = -5158.5
# Row 140
$ws.Range("H140").Value = 4167.909
$ws.Range("I140").Value = 2641.5
$ws.Range("K140").Value = 7924.5
$ws.Range("M140").Value = -2744.5

$ws = $wb.Worksheets.Item("GSM")
# Row 21
$ws.Range("H21").Value = 81666.664
$ws.Range("J21").Value = 81666.664
$ws.Range("L21").Value = 81666.664
$ws.Range("N21").Value = -82012.664
# Row 30
$ws.Range("H30").Value = 81666.664
$ws.Range("J30").Value = 81666.664
$ws.Range("L30").Value = 81666.664
$ws.Range("N30").Value = -81876.664
# Row 58
$ws.Range("H58").Value = 24499.5
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 24499.5
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 24499.5
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -25053.5

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 3133
$ws.Range("I68").Value = 3133
$ws.Range("K68").Value = 3133
$ws.Range("M68").Value = -2384
# Row 71
$ws.Range("H71").Value = 3133
$ws.Range("I71").Value = 3133
$ws.Range("K71").Value = 15665
$ws.Range("M71").Value = -11921
# Row 132
$ws.Range("H132").Value = 2830.83
$ws.Range("I132").Value = 1991.4375
$ws.Range("J132").Value = 4109.905
$ws.Range("K132").Value = 5974.3125
$ws.Range("L132").Value = 12329.715
$ws.Range("M132").Value = -3444.3125
$ws.Range("N132").Value = -17389.715
# Row 136
$ws.Range("H136").Value = 3645.55
$ws.Range("I136").Value = 2686.6155
$ws.Range("J136").Value = 5426.4287
$ws.Range("K136").Value = 8059.8465
$ws.Range("L136").Value = 16279.2861
$ws.Range("M136").Value = -5509.8465
$ws.Range("N136").Value = -21379.2861

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 1307.7727
$ws.Range("I122").Value = 1228
$ws.Range("J122").Value = 1579
$ws.Range("K122").Value = 3684
$ws.Range("L122").Value = 4737
$ws.Range("M122").Value = -1234
$ws.Range("N122").Value = -9637
# Row 126
$ws.Range("H126").Value = 4797.231
$ws.Range("I126").Value = 2821
$ws.Range("J126").Value = 9243.75
$ws.Range("K126").Value = 8463
$ws.Range("L126").Value = 27731.25
$ws.Range("M126").Value = -5993
$ws.Range("N126").Value = -32671.25
# Row 132
$ws.Range("H132").Value = 3755.4092
$ws.Range("I132").Value = 3827.628
$ws.Range("K132").Value = 11482.884
$ws.Range("M132").Value = -8952.884
# Row 135
$ws.Range("H135").Value = 57598
$ws.Range("J135").Value = 57598
$ws.Range("L135").Value = 57598
$ws.Range("N135").Value = -67738
# Row 136
$ws.Range("H136").Value = 679.3
$ws.Range("I136").Value = 532.55554
$ws.Range("K136").Value = 1597.66662
$ws.Range("M136").Value = 952.33338
